# Update "想去人数" (wanted-to-go count) figures on the two sheets that list
# the full "展览" dataset: 展览 (Exhibitions) and 全部类型 (All types).
#
# 展览 sheet layout: row5 = AP动漫游戏嘉年华, row6 = 布谷鸟动漫展4th,
#                    row9 = AB动漫游戏嘉年华
# 全部类型 sheet layout: row5 = AP动漫游戏嘉年华, row6 = 布谷鸟动漫展4th,
#                        row10 = AB动漫游戏嘉年华 (shifted by the extra
#                        "浪漫古典" concert row present on this sheet)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 2670
$ws1.Range("F6").Value = 1891
$ws1.Range("F9").Value = 941

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 2670
$ws4.Range("F6").Value = 1891
$ws4.Range("F10").Value = 941
